# Adds CDWA, CDWA Lite and CIDOC/CRM columns (G, H, I) to the metadata
# alignment table, plus a new "kept_at" row (row 11) for current location.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (row 1) ---
$ws.Range("G1").Value = "CDWA"
$ws.Range("H1").Value = "CDWA Lite"
$ws.Range("I1").Value = "CIDOC/CRM"

# --- created_by (row 2) ---
$ws.Range("G2").Value = "4.1.3. Creator Identity"
$ws.Range("H2").Value = "<cdwalite:nameCreator>"
$ws.Range("I2").Value = "<crm:P94_was_created by>"

# --- has_contributor (row 3) ---
$ws.Range("G3").Value = "4.1.3. Creator Identity"
$ws.Range("H3").Value = "<cdwalite:nameCreator>"
$ws.Range("I3").Value = "<crm:P11_had_participant>"

# --- published_by (row 4) ---
$ws.Range("G4").Value = "/"
$ws.Range("H4").Value = "/"
$ws.Range("I4").Value = "/"

# --- has_place_of_creation (row 5) ---
$ws.Range("G5").Value = "4.3. Creation Place/Original Location"
$ws.Range("H5").Value = "/"
$ws.Range("I5").Value = "/"

# --- created_in (row 6) ---
$ws.Range("G6").Value = "4.2. Creation Date"
$ws.Range("H6").Value = "<cdwalite:earliestDate>, <cdwalite:latestDate>"
$ws.Range("I6").Value = "<crm:P4_has_time-span>"

# --- has_type (row 7) ---
$ws.Range("G7").Value = "2.1. Classification Term"
$ws.Range("H7").Value = "<cdwalite:classification>"
$ws.Range("I7").Value = "<crm:P2_has_type>"

# --- subject (row 8) ---
$ws.Range("G8").Value = "16. SUBJECT MATTER"
$ws.Range("H8").Value = "<cdwalite:subjectTerm>"
$ws.Range("I8").Value = "<crm:P129_is_about>, <crm:P62_depicts>"

# --- relates_to (row 9) ---
$ws.Range("G9").Value = "/"
$ws.Range("H9").Value = "/"
$ws.Range("I9").Value = "<crm:P67_refers_to>"

# --- issued_in (row 10) ---
$ws.Range("G10").Value = "/"
$ws.Range("H10").Value = "/"
$ws.Range("I10").Value = "<crm:P148_is_component_of>"

# --- kept_at (new row 11) ---
$ws.Range("A11").Value = "kept_at"
$ws.Range("B11").Value = "<dc:coverage>"
$ws.Range("G11").Value = "21. CURRENT LOCATION"
$ws.Range("H11").Value = "<cdwalite:locationName>"
$ws.Range("I11").Value = "<crm:P55_has_current_location>"

# --- Column widths (auto-fit-like explicit widths, matching the author's
#     final layout after widening columns to fit the new content) ---
$ws.Columns.Item(1).ColumnWidth = 21.140625
$ws.Columns.Item(2).ColumnWidth = 15.7109375
$ws.Columns.Item(3).ColumnWidth = 50.85546875
$ws.Columns.Item(4).ColumnWidth = 131.28515625
$ws.Columns.Item(5).ColumnWidth = 34
$ws.Columns.Item(6).ColumnWidth = 30.140625
$ws.Columns.Item(7).ColumnWidth = 33.85546875
$ws.Columns.Item(8).ColumnWidth = 43.42578125
$ws.Columns.Item(9).ColumnWidth = 38

# --- Row 8 shrinks back to a single-line height now the row no longer needs
#     as much vertical room ---
$ws.Rows.Item(8).RowHeight = 30

# --- View state: scrolled right one column, with H10 selected ---
$ws.Range("H10").Select()
$excel.ActiveWindow.ScrollColumn = 2
